# Update "Inscritos", "Pagos" and "Inscrições homologadas" figures in the
# Resumo de Inscricoes worksheet (data refresh from the source system).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E5").Value = 31
$ws.Range("E8").Value = 47
$ws.Range("E14").Value = 41
$ws.Range("E15").Value = 108
$ws.Range("E16").Value = 330
$ws.Range("F16").Value = 95
$ws.Range("H16").Value = 183

$wb.Save()
